{"js": "// The template placeholder \"{{ photo | image(width=cm(4)) }}\" should read\n// \"{{ photo | photo(width=cm(4)) }}\" \u2014 i.e. the Jinja filter name \"image\"\n// is renamed to \"photo\" (commit: \"Changed image to photo\").\n//\n// Search the document body for the literal run of text \"image\" (it is\n// unique in this document \u2014 it only occurs inside that placeholder) and\n// replace it in place with \"photo\".\nconst body = context.document.body;\nconst results = body.search(\"image\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"photo\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The template placeholder \"{{ photo | image(width=cm(4)) }}\" should read\n# \"{{ photo | photo(width=cm(4)) }}\" \u2014 i.e. the Jinja filter name \"image\"\n# is renamed to \"photo\" (commit: \"Changed image to photo\").\n#\n# \"image\" is unique in the document (it only occurs inside that one\n# placeholder), so a simple Find/Replace over the whole document body is\n# sufficient and precise.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"image\"\n$find.Replacement.Text = \"photo\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n# wdReplaceAll = 2\n$find.Execute([ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, [ref]$null, 2) | Out-Null\n"}
